$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the expected-result string for case01 (remove stray "1" after access_token)
$ws.Range("N2").Value = "access_token,expires_in"

# Fix the JSON path expression to use the JSONPath "$" root prefix
$ws.Range("L3").Value = '$.access_token'

# Update the expected regex match value for case02's tag name
$ws.Range("N4").Value = '{"tag":{"id":(.+?),"name":"sss"}}'

# Reflect the final view/selection state on the sheet
$ws.Range("N4").Select()
